$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text value, and whether the text
# looks numeric (needs a leading apostrophe + style reset so Excel
# keeps it as literal text instead of silently coercing it to a
# Number and dropping significant trailing/leading zeros).
$updates = @(
    @{Cell='E2'; Value='  +0.72%  '; Numeric=$false}
    @{Cell='D3'; Value='1.878.61'; Numeric=$false}
    @{Cell='E3'; Value='  +1.12%  '; Numeric=$false}
    @{Cell='D4'; Value='0.9999'; Numeric=$true}
    @{Cell='E4'; Value='  -0.05%  '; Numeric=$false}
    @{Cell='D5'; Value='0.7202'; Numeric=$true}
    @{Cell='E5'; Value='  +1.64%  '; Numeric=$false}
    @{Cell='D6'; Value='240.35'; Numeric=$true}
    @{Cell='E6'; Value='  +0.81%  '; Numeric=$false}
    @{Cell='D7'; Value='0.9999'; Numeric=$true}
    @{Cell='E7'; Value='  -0.07%  '; Numeric=$false}
    @{Cell='D8'; Value='0.07814'; Numeric=$true}
    @{Cell='E8'; Value='  -2.06%  '; Numeric=$false}
    @{Cell='D9'; Value='0.3119'; Numeric=$true}
    @{Cell='E9'; Value='  +2.99%  '; Numeric=$false}
    @{Cell='D10'; Value='25.05'; Numeric=$true}
    @{Cell='E10'; Value='  +6.84%  '; Numeric=$false}
    @{Cell='D11'; Value='0.08251'; Numeric=$true}
    @{Cell='E11'; Value='  +0.61%  '; Numeric=$false}
    @{Cell='D12'; Value='1.877.13'; Numeric=$false}
    @{Cell='E12'; Value='  +4.43%  '; Numeric=$false}
    @{Cell='D13'; Value='0.7288'; Numeric=$true}
    @{Cell='E13'; Value='  +3.69%  '; Numeric=$false}
    @{Cell='D14'; Value='5.298'; Numeric=$true}
    @{Cell='E14'; Value='  +2.32%  '; Numeric=$false}
    @{Cell='D15'; Value='91.44'; Numeric=$true}
    @{Cell='E15'; Value='  +2.04%  '; Numeric=$false}
    @{Cell='D16'; Value='29.493.91'; Numeric=$false}
    @{Cell='E16'; Value='  +1.50%  '; Numeric=$false}
    @{Cell='D17'; Value='5.947'; Numeric=$true}
    @{Cell='E17'; Value='  +2.40%  '; Numeric=$false}
    @{Cell='D18'; Value='247.00'; Numeric=$true}
    @{Cell='E18'; Value='  +3.94%  '; Numeric=$false}
    @{Cell='D19'; Value='0.000007882'; Numeric=$true}
    @{Cell='E19'; Value='  -0.04%  '; Numeric=$false}
    @{Cell='E20'; Value='  +0.64%  '; Numeric=$false}
    @{Cell='D21'; Value='0.9990'; Numeric=$true}
    @{Cell='E21'; Value='  -0.06%  '; Numeric=$false}
    @{Cell='D22'; Value='7.972'; Numeric=$true}
    @{Cell='E22'; Value='  +6.74%  '; Numeric=$false}
    @{Cell='D23'; Value='1.000'; Numeric=$true}
    @{Cell='E23'; Value='  -0.07%  '; Numeric=$false}
    @{Cell='D24'; Value='0.1572'; Numeric=$true}
    @{Cell='E24'; Value='  +9.37%  '; Numeric=$false}
    @{Cell='D25'; Value='164.06'; Numeric=$true}
    @{Cell='E25'; Value='  +0.72%  '; Numeric=$false}
    @{Cell='D26'; Value='9.052'; Numeric=$true}
    @{Cell='E26'; Value='  +1.68%  '; Numeric=$false}
    @{Cell='D27'; Value='18.33'; Numeric=$true}
    @{Cell='E27'; Value='  +1.29%  '; Numeric=$false}
    @{Cell='D28'; Value='1.368'; Numeric=$true}
    @{Cell='E28'; Value='  -3.80%  '; Numeric=$false}
    @{Cell='D29'; Value='1.488'; Numeric=$true}
    @{Cell='E29'; Value='  +0.56%  '; Numeric=$false}
    @{Cell='D30'; Value='4.387'; Numeric=$true}
    @{Cell='E30'; Value='  +0.61%  '; Numeric=$false}
    @{Cell='D31'; Value='4.153'; Numeric=$true}
    @{Cell='E31'; Value='  +3.32%  '; Numeric=$false}
    @{Cell='D32'; Value='0.05284'; Numeric=$true}
    @{Cell='E32'; Value='  +1.88%  '; Numeric=$false}
    @{Cell='D33'; Value='1.947'; Numeric=$true}
    @{Cell='E33'; Value='  +1.60%  '; Numeric=$false}
    @{Cell='E34'; Value='  +3.84%  '; Numeric=$false}
    @{Cell='E35'; Value='  +1.55%  '; Numeric=$false}
    @{Cell='D36'; Value='2.678'; Numeric=$true}
    @{Cell='E36'; Value='  +0.27%  '; Numeric=$false}
    @{Cell='E37'; Value='  +0.79%  '; Numeric=$false}
    @{Cell='D38'; Value='1.236.43'; Numeric=$false}
    @{Cell='D39'; Value='2.723'; Numeric=$true}
    @{Cell='E39'; Value='  +0.08%  '; Numeric=$false}
    @{Cell='D40'; Value='0.9085'; Numeric=$true}
    @{Cell='E40'; Value='  -2.47%  '; Numeric=$false}
    @{Cell='D41'; Value='73.78'; Numeric=$true}
    @{Cell='E41'; Value='  +5.18%  '; Numeric=$false}
    @{Cell='D42'; Value='6.105'; Numeric=$true}
    @{Cell='E42'; Value='  +3.58%  '; Numeric=$false}
    @{Cell='E43'; Value='  -0.03%  '; Numeric=$false}
    @{Cell='D44'; Value='103.76'; Numeric=$true}
    @{Cell='E44'; Value='  +1.38%  '; Numeric=$false}
    @{Cell='D45'; Value='0.5346'; Numeric=$true}
    @{Cell='E45'; Value='  +0.31%  '; Numeric=$false}
    @{Cell='B46'; Value='SynthetixNetwork'; Numeric=$false}
    @{Cell='C46'; Value='https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'; Numeric=$false}
    @{Cell='D46'; Value='2.928'; Numeric=$true}
    @{Cell='E46'; Value='  +13.10%  '; Numeric=$false}
    @{Cell='B47'; Value='RenderToken'; Numeric=$false}
    @{Cell='C47'; Value='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; Numeric=$false}
    @{Cell='D47'; Value='1.764'; Numeric=$true}
    @{Cell='E47'; Value='  +0.18%  '; Numeric=$false}
    @{Cell='B48'; Value='BabyDogeCoin'; Numeric=$false}
    @{Cell='C48'; Value='https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'; Numeric=$false}
    @{Cell='D48'; Value='0.00000000120'; Numeric=$true}
    @{Cell='E48'; Value='  +0.46%  '; Numeric=$false}
    @{Cell='B49'; Value='EnergySwap'; Numeric=$false}
    @{Cell='C49'; Value='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; Numeric=$false}
    @{Cell='D49'; Value='9.297'; Numeric=$true}
    @{Cell='E49'; Value='  +1.26%  '; Numeric=$false}
    @{Cell='B50'; Value='TheSandbox'; Numeric=$false}
    @{Cell='C50'; Value='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; Numeric=$false}
    @{Cell='D50'; Value='0.4332'; Numeric=$true}
    @{Cell='E50'; Value='  +1.90%  '; Numeric=$false}
    @{Cell='D51'; Value='7.087'; Numeric=$true}
    @{Cell='E51'; Value='  +2.00%  '; Numeric=$false}
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.Numeric) {
        # Leading apostrophe forces Excel to store the literal text
        # (prevents "247.00" -> 247, "0.9990" -> 0.999, etc.).
        $range.Value = [string]::Concat("'", $u.Value)
        # Drop the quote-prefix formatting Excel applied so the cell
        # style matches a plain unstyled text cell again.
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Value
    }
}
